$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 ("Sound"): fix the "volune" typo and add the two new fields.
# Resulting shared strings are appended in the order muted, volume, isPlaying,
# so write the cells in the matching order (E, D, then B); C14 ("melody") is
# rewritten too so it keeps its place in the table.
$ws.Range("E14").Value = "muted"
$ws.Range("D14").Value = "volume"
$ws.Range("C14").Value = "melody"
$ws.Range("B14").Value = "isPlaying"

# Move/keep the active selection on D14, as in the saved workbook.
$ws.Range("D14").Select()
